# Updated symbol list on Wed Jan 18 13:45:56 UTC 2023 with GitHub Actions
#
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto symbols whose quotes moved. Values are written with a leading
# apostrophe so Excel keeps them as literal text (matching the sheet's
# existing inline-string cells, e.g. "301.72" / "0.05%") instead of
# auto-converting them into numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value  = "'301.72"
$ws.Range("E2").Value  = "'0.05%"

$ws.Range("D3").Value  = "'32.41"
$ws.Range("E3").Value  = "'1.64%"

$ws.Range("D4").Value  = "'5.061"
$ws.Range("E4").Value  = "'-1.13%"

$ws.Range("D5").Value  = "'0.07674"
$ws.Range("E5").Value  = "'-2.13%"

$ws.Range("D6").Value  = "'2.131"
$ws.Range("E6").Value  = "'-5.12%"

$ws.Range("D7").Value  = "'7.846"
$ws.Range("E7").Value  = "'0.45%"

$ws.Range("D8").Value  = "'0.9202"
$ws.Range("E8").Value  = "'-0.82%"

$ws.Range("D9").Value  = "'0.1756"
$ws.Range("E9").Value  = "'-0.87%"

$ws.Range("D10").Value = "'0.07933"
$ws.Range("E10").Value = "'3.66%"

$ws.Range("D11").Value = "'0.08419"
$ws.Range("E11").Value = "'-4.83%"

$ws.Range("D12").Value = "'0.03063"
$ws.Range("E12").Value = "'-1.22%"

$ws.Range("D13").Value = "'0.09982"
$ws.Range("E13").Value = "'-0.40%"

$ws.Range("D14").Value = "'0.001512"
$ws.Range("E14").Value = "'-0.29%"

$ws.Range("D15").Value = "'0.005726"
$ws.Range("E15").Value = "'-2.37%"

$ws.Range("D17").Value = "'3.463"
$ws.Range("E17").Value = "'0.02%"

$ws.Range("D18").Value = "'3.778"
$ws.Range("E18").Value = "'-0.82%"

$ws.Range("E19").Value = "'-4.40%"

$ws.Range("D20").Value = "'0.3344"
$ws.Range("E20").Value = "'1.58%"

$ws.Range("D21").Value = "'0.1318"
$ws.Range("E21").Value = "'-0.72%"

$ws.Range("D22").Value = "'4.284"
$ws.Range("E22").Value = "'-0.91%"

$ws.Range("D23").Value = "'0.1978"
$ws.Range("E23").Value = "'10.42%"

$ws.Range("D24").Value = "'0.04543"
$ws.Range("E24").Value = "'-1.31%"

$ws.Range("D25").Value = "'0.001235"
$ws.Range("E25").Value = "'-1.31%"

$ws.Range("D26").Value = "'0.004827"
$ws.Range("E26").Value = "'7.65%"

$ws.Range("D27").Value = "'0.0001249"
$ws.Range("E27").Value = "'-0.04%"

$ws.Range("D39").Value = "'0.01714"
$ws.Range("E39").Value = "'-3.79%"

$ws.Range("D40").Value = "'0.04677"
$ws.Range("E40").Value = "'-2.19%"

$ws.Range("D41").Value = "'0.007549"
$ws.Range("E41").Value = "'2.59%"

$ws.Range("D42").Value = "'0.1355"
$ws.Range("E42").Value = "'-0.60%"

$ws.Range("D43").Value = "'0.002329"
$ws.Range("E43").Value = "'6.35%"

$ws.Range("D44").Value = "'0.01061"
$ws.Range("E44").Value = "'7.74%"

$ws.Range("D45").Value = "'0.00006199"
$ws.Range("E45").Value = "'-1.04%"

$ws.Range("E46").Value = "'-0.03%"

$ws.Range("D47").Value = "'1.266"
$ws.Range("E47").Value = "'80.70%"

$ws.Range("D48").Value = "'0.002998"
$ws.Range("E48").Value = "'-62.47%"

$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.03%"

$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.03%"
